# Generate Report for Handback
# a.md has been handed back (in sync with en-US); update the Overview sheet and the
# per-locale (zh-cn / de-de) detail sheets to reflect the new handback status,
# handback datetime/name, and clear the stale "not latest" error for a.md.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: row 2 is a.md. Status columns (zh-cn / de-de) move from
# "In Translation" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# Widen the status columns so the longer text is readable.
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------------
# zh-cn sheet: row 2 is a.md.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("L2").Value = "2017-02-17 07:16:41"
$wsZhCn.Range("M2").Value = "TestHandback_201702170316"
$wsZhCn.Range("R2").Value = ""

$wsZhCn.Columns.Item(3).ColumnWidth = 29.14
$wsZhCn.Columns.Item(13).ColumnWidth = 27.14

# ---------------------------------------------------------------------------
# de-de sheet: row 2 is a.md.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("L2").Value = "2017-02-17 07:17:04"
$wsDeDe.Range("M2").Value = "TestHandback_201702170316"
$wsDeDe.Range("R2").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = 29.14
$wsDeDe.Columns.Item(13).ColumnWidth = 27.14
